$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B:F per row (2-11), and new N values for G per row (2-11)
$values = @{
    2  = @(0.2736504459652551, 0.6164471910808867, 0.5099410358193527, 0.7141015584770507, 0.6776623592093396, 19)
    3  = @(0.1679866110975223, 0.4776555476168131, 0.3912275594536091, 0.6254818618102439, 0.6199689291560088, 18)
    4  = @(0.2407374606666166, 0.3283263025081246, 0.1549134947928615, 0.393590516644979,  0.320965582950402,  17)
    5  = @(0.3026896478973899, 0.3981352914125425, 0.2125958988159086, 0.4610812280020827, 0.3592211959175575, 16)
    6  = @(0.3313557959155047, 0.3928030709233679, 0.2145494597484939, 0.463194839941567,  0.3350151322728833, 15)
    7  = @(0.3718062059916627, 0.3799713553340741, 0.2141292450220557, 0.4627410129025259, 0.2858794074584022, 14)
    8  = @(0.4009368084426636, 0.4127405556116039, 0.227594880374374,  0.4770690519981086, 0.2691002335147636, 13)
    9  = @(0.4356636671697872, 0.4368058703952601, 0.2362224927271332, 0.4860272551278716, 0.2250325111912971, 12)
    10 = @(0.4167250561910649, 0.4167250561910649, 0.2156249438995315, 0.4643543301182099, 0.2148527137047464, 11)
    11 = @(0.3832484536769914, 0.3832484536769914, 0.1811246737723569, 0.4255874455060403, 0.1950649365568976, 10)
}

foreach ($row in $values.Keys) {
    $vals = $values[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
